$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Column D (Price) values are text that can look numeric (e.g. "1.0000"),
# so we force the Text number format before assigning, then restore the
# "Normal" cell style afterwards so no stray formatting is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.204.39'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.841.77'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.68'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07458'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2955'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.27'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07746'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.837.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.017'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6732'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.175'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008734'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.202.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.092.90'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '226.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.193'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.622'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1393'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.514'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.138'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.210'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.043'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05385'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7469'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.158'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.64%  '
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.300.10'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01797'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.354'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9044'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.08323'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.991.14'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000123'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '65.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.5141'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.753'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.94%  '
